$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1945.5454
$ws.Range("I6").Value = 125.25
$ws.Range("K6").Value = 375.75
$ws.Range("M6").Value = -263.75

$ws.Range("H40").Value = 1636.55
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350

$ws.Range("H96").Value = 406.66666
$ws.Range("I96").Value = 210
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 630
$ws.Range("L96").Value = 2400
$ws.Range("M96").Value = 743
$ws.Range("N96").Value = -5146

$ws.Range("H99").Value = 1533.7142
$ws.Range("I99").Value = 311.5
$ws.Range("J99").Value = 3163.3333
$ws.Range("K99").Value = 934.5
$ws.Range("L99").Value = 9489.999899999999
$ws.Range("M99").Value = 563.5
$ws.Range("N99").Value = -12485.9999

$ws.Range("H101").Value = 5600
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 5600
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 16800
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -20044

$ws.Range("H116").Value = 3338.75
$ws.Range("J116").Value = 4396.25
$ws.Range("L116").Value = 4396.25
$ws.Range("N116").Value = -11280.25

$ws.Range("H129").Value = 1200.8077
$ws.Range("I129").Value = 2199
$ws.Range("J129").Value = 1070.6086
$ws.Range("K129").Value = 6597
$ws.Range("L129").Value = 3211.8258
$ws.Range("M129").Value = -1597
$ws.Range("N129").Value = -13211.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5760.67
$ws.Range("I32").Value = 5350.086
$ws.Range("J32").Value = 11215.571
$ws.Range("K32").Value = 5350.086
$ws.Range("L32").Value = 11215.571
$ws.Range("M32").Value = -5063.086
$ws.Range("N32").Value = -11789.571

$ws.Range("H61").Value = 2053.7778
$ws.Range("I61").Value = 1135.4286
$ws.Range("J61").Value = 3042.7693
$ws.Range("K61").Value = 1135.4286
$ws.Range("L61").Value = 3042.7693
$ws.Range("M61").Value = -923.4286
$ws.Range("N61").Value = -3466.7693

$ws.Range("H74").Value = 1218.8
$ws.Range("I74").Value = 1231.2222
$ws.Range("J74").Value = 1107
$ws.Range("K74").Value = 1231.2222
$ws.Range("L74").Value = 1107
$ws.Range("M74").Value = -357.2221999999999
$ws.Range("N74").Value = -2855

$ws.Range("H77").Value = 1218.8
$ws.Range("I77").Value = 1231.2222
$ws.Range("J77").Value = 1107
$ws.Range("K77").Value = 6156.111
$ws.Range("L77").Value = 5535
$ws.Range("M77").Value = -1788.111
$ws.Range("N77").Value = -14271

$ws.Range("H102").Value = 4744.875
$ws.Range("I102").Value = 3994.1428
$ws.Range("K102").Value = 3994.1428
$ws.Range("M102").Value = -2372.1428

$ws.Range("H122").Value = 2382.7026
$ws.Range("I122").Value = 1866.5186
$ws.Range("K122").Value = 5599.5558
$ws.Range("M122").Value = -3149.5558

$ws.Range("H132").Value = 2227.4565
$ws.Range("I132").Value = 1835.4572
$ws.Range("K132").Value = 5506.3716
$ws.Range("M132").Value = -2976.3716

$ws.Range("H136").Value = 2053.7778
$ws.Range("I136").Value = 1135.4286
$ws.Range("J136").Value = 3042.7693
$ws.Range("K136").Value = 3406.2858
$ws.Range("L136").Value = 9128.3079
$ws.Range("M136").Value = -856.2857999999997
$ws.Range("N136").Value = -14228.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 565
$ws.Range("I5").Value = 279
$ws.Range("J5").Value = 922.5
$ws.Range("K5").Value = 279
$ws.Range("L5").Value = 922.5
$ws.Range("M5").Value = -166
$ws.Range("N5").Value = -1148.5

$ws.Range("H94").Value = 813.1667
$ws.Range("I94").Value = 816
$ws.Range("J94").Value = 799
$ws.Range("K94").Value = 816
$ws.Range("L94").Value = 799
$ws.Range("M94").Value = -365
$ws.Range("N94").Value = -1701

$ws.Range("H134").Value = 2156.4866
$ws.Range("I134").Value = 1540.5518
$ws.Range("J134").Value = 4389.25
$ws.Range("K134").Value = 4621.6554
$ws.Range("L134").Value = 13167.75
$ws.Range("M134").Value = -2086.6554
$ws.Range("N134").Value = -18237.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10641129
$ws.Range("I58").Value = 1919.258
$ws.Range("J58").Value = 31254600
$ws.Range("K58").Value = 1919.258
$ws.Range("L58").Value = 31254600
$ws.Range("M58").Value = -1716.258
$ws.Range("N58").Value = -31255006

$ws.Range("H94").Value = 5592
$ws.Range("I94").Value = 6202
$ws.Range("K94").Value = 6202
$ws.Range("M94").Value = -5751

$ws.Range("H107").Value = 1332.4348
$ws.Range("I107").Value = 491.65
$ws.Range("J107").Value = 6937.6665
$ws.Range("K107").Value = 491.65
$ws.Range("L107").Value = 6937.6665
$ws.Range("M107").Value = 1428.35
$ws.Range("N107").Value = -10777.6665

$ws.Range("H134").Value = 4230
$ws.Range("I134").Value = 650
$ws.Range("J134").Value = 5125
$ws.Range("K134").Value = 1950
$ws.Range("L134").Value = 15375
$ws.Range("M134").Value = 585
$ws.Range("N134").Value = -20445

$ws.Range("H136").Value = 10641129
$ws.Range("I136").Value = 1919.258
$ws.Range("J136").Value = 31254600
$ws.Range("K136").Value = 5757.774
$ws.Range("L136").Value = 93763800
$ws.Range("M136").Value = -3207.774
$ws.Range("N136").Value = -93768900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3444
$ws.Range("I119").Value = 2008
$ws.Range("J119").Value = 3982.5
$ws.Range("K119").Value = 6024
$ws.Range("L119").Value = 11947.5
$ws.Range("M119").Value = -1186
$ws.Range("N119").Value = -21623.5

$ws.Range("H140").Value = 1686.1154
$ws.Range("I140").Value = 842.55
$ws.Range("J140").Value = 4498
$ws.Range("K140").Value = 2527.65
$ws.Range("L140").Value = 13494
$ws.Range("M140").Value = 2652.35
$ws.Range("N140").Value = -23854

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3983.9312
$ws.Range("I132").Value = 4586.3335
$ws.Range("J132").Value = 3558.7058
$ws.Range("K132").Value = 13759.0005
$ws.Range("L132").Value = 10676.1174
$ws.Range("M132").Value = -11229.0005
$ws.Range("N132").Value = -15736.1174

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H46").Value = 1199.8182
$ws.Range("I46").Value = 859.4386
$ws.Range("J46").Value = 3355.5557
$ws.Range("K46").Value = 859.4386
$ws.Range("L46").Value = 3355.5557
$ws.Range("M46").Value = -671.4386
$ws.Range("N46").Value = -3731.5557

$ws.Range("H132").Value = 1687.8679
$ws.Range("I132").Value = 1099.0444
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3297.1332
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -767.1332000000002
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2280
$ws.Range("I96").Value = 1800
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 1800
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -427
$ws.Range("N96").Value = -5746

$ws.Range("H100").Value = 523.75
$ws.Range("I100").Value = 450
$ws.Range("J100").Value = 597.5
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 1195
$ws.Range("M100").Value = -359
$ws.Range("N100").Value = -2277

$ws.Range("H132").Value = 3340.5967
$ws.Range("I132").Value = 1316.9259
$ws.Range("K132").Value = 3950.7777
$ws.Range("M132").Value = -1420.7777

$ws.Range("H136").Value = 1197.05
$ws.Range("I136").Value = 599.37933
$ws.Range("J136").Value = 2772.7273
$ws.Range("K136").Value = 1798.13799
$ws.Range("L136").Value = 8318.1819
$ws.Range("M136").Value = 751.8620100000001
$ws.Range("N136").Value = -13418.1819
